# Generate Report for Handoff
#
# Marks the six "Ready for handoff" rows (7,8,9,10,12,13 -- row 11 is a
# duplicate-content row that is intentionally skipped) as handoff type "ht"
# on both language sheets, and refreshes the handoff-generation timestamps
# that accompany that run: the zh-cn sheet's "Latest Handoff Datetime"
# column, the de-de sheet's "Latest Handoff Datetime" column, and the
# Overview sheet's "Latest HO Xliff Generate Date" column.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 12, 13)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value2 = "ht"
    $wsZhCn.Range("H$r").Value2 = "2016-08-31 00:23:33"
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value2 = "ht"
    $wsDeDe.Range("H$r").Value2 = "2016-08-31 00:23:38"
}

$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value2 = "2016-08-31 00:23:38"
}
